$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Split the R12/R13 "10k" resistor row into two rows -----------------
# Before: row 28 = Qty 19, Parts "R2, ..., R13, ... R36" (10k / CRCW040210K0FKED)
# After : row 28 = Qty 18, Parts without R13; new row 29 = R13 (10k 0.1% / CPF0402B10KE1)

# Insert a new blank row right below the "10k" row (row 28), pushing
# everything from row 29 down (this also shifts the merged header rows and
# the section rows below automatically, matching native Excel behaviour).
$ws.Rows.Item(29).Insert()

# Update the existing 10k row (row 28): drop the qty by one and remove
# "R13" from the comma-separated Parts list.
$ws.Cells.Item(28, 1).Value2 = 18
$ws.Cells.Item(28, 4).Value2 = "R2, R3, R4, R5, R6, R7, R8, R10, R11, R15, R17, R21, R22, R23, R32, R33, R34, R36"

# Fill in the new row 29 for R13 (10k 0.1%).
$ws.Cells.Item(29, 1).Value2 = 1
$ws.Cells.Item(29, 2).Value2 = "10k 0.1%"
$ws.Cells.Item(29, 3).Value2 = "R0402"
$ws.Cells.Item(29, 4).Value2 = "R13"
$ws.Cells.Item(29, 5).Value2 = "CPF0402B10KE1"
$ws.Cells.Item(29, 6).Value2 = "SMD"
$ws.Cells.Item(29, 7).Value2 = "-"

# --- Update the R12 row (now row 31) with the new 0.1% value/MPN --------
$ws.Cells.Item(31, 2).Value2 = "34k 0.1%"
$ws.Cells.Item(31, 5).Value2 = "ERA2AEB3402X"

# --- Fix up the stale hidden _FilterDatabase defined name ----------------
# It previously pointed at $A$43:$F$49; after inserting a row above it, it
# should now point at $A$44:$F$50.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='harp expander v1'!`$A`$44:`$F`$50"
    }
}

# --- Keep the view selection in a sensible spot (cosmetic) ---------------
$ws.Application.Goto($ws.Range("A7"), $true)
$ws.Range("D24").Select()
